$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; this shifts the existing rows 28-149 down to 29-150
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new data record
$ws.Cells.Item(28, 1).Value = 4
$ws.Cells.Item(28, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(28, 3).Value = "Los Lagos"
$ws.Cells.Item(28, 4).Value = 44575
$ws.Cells.Item(28, 5).Value = 10
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100108
$ws.Cells.Item(28, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(28, 9).Value = 100108002
$ws.Cells.Item(28, 10).Value = "Mango"
$ws.Cells.Item(28, 11).Value = "Sin especificar"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 200
$ws.Cells.Item(28, 14).Value = 8000
$ws.Cells.Item(28, 15).Value = 8500
$ws.Cells.Item(28, 16).Value = 8250
$ws.Cells.Item(28, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(28, 18).Value = "Perú"
$ws.Cells.Item(28, 19).Value = 2062
$ws.Cells.Item(28, 20).Value = 4

# Make sure the date cell keeps the same date/time number format as the rest of column D
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat()
